# Applies the two content edits described by the commit diff:
#   1. Inside the paragraph about "en la actualidad se producen grandes
#      cantidades de libros..." the phrase "desechados a lo que conlleva"
#      becomes "desechados, lo que conlleva" (the " a " between
#      "desechados" and "lo" is replaced by ", ").
#   2. The whole paragraph that reads "Valoración como encuesta" (an
#      underlined heading-like line right after the "Las fuentes de
#      información..." paragraph) is removed completely, paragraph mark
#      and all, while the empty paragraph that follows it is kept.

$d = $word.ActiveDocument

# --- 1) "desechados a lo" -> "desechados, lo" -------------------------
$find = $d.Content
$found = $find.Find.Execute(
    "desechados a lo que conlleva",  # find
    $true,                           # MatchCase
    $false,                          # MatchWholeWord
    $false,                          # MatchWildcards
    $false,                          # MatchSoundsLike
    $false,                          # MatchAllWordForms
    $true,                           # Forward
    1,                               # Wrap (wdFindContinue)
    $false,                          # Format
    "desechados, lo que conlleva",   # Replace
    2                                # Replace (wdReplaceAll)
)
if (-not $found) {
    throw "Could not find the 'desechados a lo que conlleva' text to update."
}

# --- 2) Remove the "Valoración como encuesta" paragraph ---------------
$target = $d.Content
$found2 = $target.Find.Execute(
    "Valoración como encuesta",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
if (-not $found2) {
    throw "Could not find the 'Valoración como encuesta' paragraph to remove."
}
# Expand the found range to the whole paragraph (wdParagraph = 4) so the
# paragraph mark is included too, then delete it outright so no empty
# paragraph is left behind in its place.
$target.Expand(4) | Out-Null
$target.Delete()
